$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.769.40'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.801.57'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.05'
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.26'
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("E8").Value = '  +0.56%  '

$ws.Range("E9").Value = '  +1.77%  '

$ws.Range("E10").Value = '  -1.08%  '

$ws.Range("E11").Value = '  +0.55%  '

$ws.Range("E12").Value = '  -0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.98'
$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.442.68'
$ws.Range("E14").Value = '  +0.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.841.83'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.54'
$ws.Range("E16").Value = '  +3.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.773.77'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.113'
$ws.Range("E19").Value = '  +0.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.67'
$ws.Range("E20").Value = '  +1.14%  '

$ws.Range("E21").Value = '  -3.08%  '

$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("E23").Value = '  +1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.27'
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("E25").Value = '  +2.12%  '

$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("E28").Value = '  +0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.945.95'
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E31").Value = '  +2.36%  '

$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("E35").Value = '  -0.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.740.74'
$ws.Range("E36").Value = '  +0.23%  '

$ws.Range("E37").Value = '  +0.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +2.93%  '

$ws.Range("E39").Value = '  +0.29%  '

$ws.Range("E40").Value = '  +0.99%  '

$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.07'
$ws.Range("E44").Value = '  +2.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.302'
$ws.Range("E45").Value = '  +1.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.87'
$ws.Range("E46").Value = '  -1.99%  '

$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '147.44'
$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +10.12%  '

$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '395.50'
$ws.Range("E50").Value = '  +1.31%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '27.09'
$ws.Range("E51").Value = '  +7.16%  '
